$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("facebook")

# Rename the sheet from "facebook" to "flipkart"
$ws.Name = "flipkart"

# Update the title cell (B1) text: 67 leading spaces + "flipcart"
$ws.Range("B1").Value = "                                                                   flipcart"
